$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'270.02"
$ws.Range("E2").Value = "'3.18%"
$ws.Range("E3").Value = "'-1.65%"
$ws.Range("D4").Value = "'4.715"
$ws.Range("E4").Value = "'0.12%"
$ws.Range("D5").Value = "'0.06104"
$ws.Range("E5").Value = "'-1.69%"
$ws.Range("D6").Value = "'6.739"
$ws.Range("E6").Value = "'0.28%"
$ws.Range("D7").Value = "'0.8568"
$ws.Range("E7").Value = "'0.79%"
$ws.Range("D8").Value = "'0.8905"
$ws.Range("E8").Value = "'-2.10%"
$ws.Range("D9").Value = "'0.1431"
$ws.Range("E9").Value = "'1.37%"
$ws.Range("D10").Value = "'0.04950"
$ws.Range("E10").Value = "'5.38%"
$ws.Range("D11").Value = "'0.07119"
$ws.Range("E11").Value = "'0.34%"
$ws.Range("D12").Value = "'0.03184"
$ws.Range("E12").Value = "'0.83%"
$ws.Range("D13").Value = "'0.09033"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("D14").Value = "'0.001540"
$ws.Range("E14").Value = "'0.27%"
$ws.Range("E15").Value = "'-1.17%"
$ws.Range("D16").Value = "'0.005945"
$ws.Range("E16").Value = "'-2.16%"
$ws.Range("E17").Value = "'-0.17%"
$ws.Range("D18").Value = "'3.173"
$ws.Range("E18").Value = "'0.05%"
$ws.Range("D19").Value = "'2.264"
$ws.Range("E19").Value = "'3.86%"
$ws.Range("D20").Value = "'0.3090"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("E21").Value = "'-0.72%"
$ws.Range("D22").Value = "'3.846"
$ws.Range("D23").Value = "'0.04242"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'-2.33%"
$ws.Range("D25").Value = "'0.004149"
$ws.Range("E25").Value = "'0.38%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("D27").Value = "'0.0001680"
$ws.Range("E27").Value = "'4.98%"
$ws.Range("D40").Value = "'0.03947"
$ws.Range("E40").Value = "'1.11%"
$ws.Range("D41").Value = "'0.1118"
$ws.Range("E41").Value = "'0.43%"
$ws.Range("D42").Value = "'0.004185"
$ws.Range("E42").Value = "'1.32%"
$ws.Range("E43").Value = "'-6.76%"
$ws.Range("E44").Value = "'-15.26%"
$ws.Range("D45").Value = "'0.00005121"
$ws.Range("E45").Value = "'-1.00%"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("E47").Value = "'-31.83%"
$ws.Range("D48").Value = "'0.9468"
$ws.Range("E48").Value = "'467.09%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.09%"
